# Update cryptocurrency price/volume figures to the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value (kept as literal text,
# matching the sheet's existing inline-string formatting for Price/Volume columns).
$updates = @{
    "D2" = "307.75"
    "E2" = "-0.49%"
    "D3" = "39.74"
    "E3" = "2.24%"
    "D4" = "5.166"
    "E4" = "1.57%"
    "D5" = "0.08109"
    "E5" = "-1.02%"
    "D6" = "1.947"
    "E6" = "-4.50%"
    "D7" = "8.152"
    "E7" = "3.07%"
    "D8" = "0.9260"
    "E8" = "-0.46%"
    "D9" = "0.1421"
    "E9" = "0.09%"
    "D10" = "0.1920"
    "E10" = "-1.59%"
    "D11" = "0.09104"
    "E11" = "-1.87%"
    "D12" = "0.03518"
    "E12" = "1.43%"
    "D13" = "0.09800"
    "E13" = "-0.45%"
    "D14" = "0.001398"
    "E14" = "-0.48%"
    "D15" = "0.005855"
    "E15" = "-0.54%"
    "D16" = "3.935"
    "E16" = "4.62%"
    "D17" = "4.237"
    "E17" = "1.41%"
    "E18" = "-2.47%"
    "E19" = "-0.71%"
    "D20" = "0.1320"
    "E20" = "1.36%"
    "D21" = "4.662"
    "E21" = "-3.48%"
    "D22" = "0.2422"
    "E23" = "-2.19%"
    "D24" = "0.001229"
    "E24" = "-0.80%"
    "D25" = "0.004362"
    "E25" = "4.50%"
    "E26" = "0.09%"
    "D27" = "0.0004001"
    "E27" = "-10.04%"
    "D39" = "0.02041"
    "E39" = "-3.91%"
    "D40" = "0.05059"
    "E40" = "-2.26%"
    "D41" = "0.007399"
    "E41" = "-0.88%"
    "D42" = "0.009796"
    "E42" = "-3.23%"
    "E43" = "-0.31%"
    "E44" = "0.10%"
    "D45" = "0.009482"
    "E45" = "-2.02%"
    "D46" = "0.00006363"
    "E46" = "0.87%"
    "D47" = "0.00000000750"
    "E47" = "-0.04%"
    "D48" = "0.002728"
    "E49" = "-18.77%"
    "D50" = "0.00002101"
    "E50" = "-0.04%"
    "D51" = "0.0002001"
    "E51" = "-0.04%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number-format before assigning so Excel keeps these
    # numeric-looking / percent-looking strings as literal text instead
    # of parsing them into numbers, then restore the default "Normal"
    # cell style so no stray formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
